# Auto-generated script: apply numeric updates to Chocobo_Profits workbook
# Each sheet corresponds to a crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns H-N hold computed market/profit figures that were refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1872.1428
$ws.Range("I28").Value = 1421
$ws.Range("J28").Value = 3000
$ws.Range("K28").Value = 1421
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = -936
$ws.Range("N28").Value = -3970
$ws.Range("H43").Value = 1910.1428
$ws.Range("I43").Value = 1083.1666
$ws.Range("J43").Value = 2530.375
$ws.Range("K43").Value = 1083.1666
$ws.Range("L43").Value = 2530.375
$ws.Range("M43").Value = -1014.1666
$ws.Range("N43").Value = -2668.375
$ws.Range("H112").Value = 1317.8948
$ws.Range("J112").Value = 1317.8948
$ws.Range("L112").Value = 3953.6844
$ws.Range("N112").Value = -6169.6844
$ws.Range("H125").Value = 3336
$ws.Range("J125").Value = 3336
$ws.Range("L125").Value = 30024
$ws.Range("N125").Value = -34944
$ws.Range("H129").Value = 849.1
$ws.Range("I129").Value = 366.66666
$ws.Range("J129").Value = 864.0205999999999
$ws.Range("K129").Value = 1099.99998
$ws.Range("L129").Value = 2592.0618
$ws.Range("M129").Value = 3900.00002
$ws.Range("N129").Value = -12592.0618
$ws.Range("H141").Value = 64877.75
$ws.Range("I141").Value = 78904.16
$ws.Range("K141").Value = 236712.48
$ws.Range("M141").Value = -231532.48

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6571.0796
$ws.Range("I32").Value = 4992.0977
$ws.Range("K32").Value = 4992.0977
$ws.Range("M32").Value = -4705.0977
$ws.Range("H61").Value = 1835.091
$ws.Range("I61").Value = 1778.25
$ws.Range("J61").Value = 1986.6666
$ws.Range("K61").Value = 1778.25
$ws.Range("L61").Value = 1986.6666
$ws.Range("M61").Value = -1566.25
$ws.Range("N61").Value = -2410.6666
$ws.Range("H122").Value = 12004.8
$ws.Range("I122").Value = 10008
$ws.Range("K122").Value = 30024
$ws.Range("M122").Value = -27574
$ws.Range("H136").Value = 1835.091
$ws.Range("I136").Value = 1778.25
$ws.Range("J136").Value = 1986.6666
$ws.Range("K136").Value = 5334.75
$ws.Range("L136").Value = 5959.9998
$ws.Range("M136").Value = -2784.75
$ws.Range("N136").Value = -11059.9998

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 4567.4
$ws.Range("I36").Value = 709.25
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 709.25
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -175.25
$ws.Range("N36").Value = -21068
$ws.Range("H115").Value = 34657.895
$ws.Range("J115").Value = 34657.895
$ws.Range("L115").Value = 34657.895
$ws.Range("N115").Value = -37791.895

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7414.64
$ws.Range("I31").Value = 1507
$ws.Range("J31").Value = 11353.066
$ws.Range("K31").Value = 1507
$ws.Range("L31").Value = 11353.066
$ws.Range("M31").Value = -1212
$ws.Range("N31").Value = -11943.066
$ws.Range("H34").Value = 7414.64
$ws.Range("I34").Value = 1507
$ws.Range("J34").Value = 11353.066
$ws.Range("K34").Value = 1507
$ws.Range("L34").Value = 11353.066
$ws.Range("M34").Value = -1305
$ws.Range("N34").Value = -11757.066
$ws.Range("H99").Value = 11115439
$ws.Range("I99").Value = 20002170
$ws.Range("J99").Value = 7025
$ws.Range("K99").Value = 20002170
$ws.Range("L99").Value = 7025
$ws.Range("M99").Value = -20000672
$ws.Range("N99").Value = -10021
$ws.Range("H107").Value = 764.6
$ws.Range("I107").Value = 627.3333
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 627.3333
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1292.6667
$ws.Range("N107").Value = -5840
$ws.Range("H115").Value = 34900
$ws.Range("J115").Value = 34900
$ws.Range("L115").Value = 34900
$ws.Range("N115").Value = -37250
$ws.Range("H126").Value = 11115439
$ws.Range("I126").Value = 20002170
$ws.Range("J126").Value = 7025
$ws.Range("K126").Value = 60006510
$ws.Range("L126").Value = 21075
$ws.Range("M126").Value = -60004040
$ws.Range("N126").Value = -26015

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 170.91667
$ws.Range("I23").Value = 175
$ws.Range("J23").Value = 170.1
$ws.Range("K23").Value = 525
$ws.Range("L23").Value = 510.3
$ws.Range("M23").Value = -290
$ws.Range("N23").Value = -980.3
$ws.Range("H134").Value = 4412.6665
$ws.Range("I134").Value = 3877.1428
$ws.Range("J134").Value = 4881.25
$ws.Range("K134").Value = 11631.4284
$ws.Range("L134").Value = 14643.75
$ws.Range("M134").Value = -6561.428400000001
$ws.Range("N134").Value = -24783.75

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5839.54
$ws.Range("I70").Value = 5514.45
$ws.Range("J70").Value = 7139.9
$ws.Range("K70").Value = 5514.45
$ws.Range("L70").Value = 7139.9
$ws.Range("M70").Value = -5244.45
$ws.Range("N70").Value = -7679.9
$ws.Range("H73").Value = 5839.54
$ws.Range("I73").Value = 5514.45
$ws.Range("J73").Value = 7139.9
$ws.Range("K73").Value = 5514.45
$ws.Range("L73").Value = 7139.9
$ws.Range("M73").Value = -4578.45
$ws.Range("N73").Value = -9011.9
$ws.Range("H132").Value = 5449.7334
$ws.Range("I132").Value = 4562.3335
$ws.Range("J132").Value = 8999.333000000001
$ws.Range("K132").Value = 13687.0005
$ws.Range("L132").Value = 26997.999
$ws.Range("M132").Value = -11157.0005
$ws.Range("N132").Value = -32057.999

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5642.143
$ws.Range("I7").Value = 2098.3333
$ws.Range("J7").Value = 8300
$ws.Range("K7").Value = 2098.3333
$ws.Range("L7").Value = 8300
$ws.Range("M7").Value = -1986.3333
$ws.Range("N7").Value = -8524
$ws.Range("H40").Value = 9000
$ws.Range("I40").Value = 7000
$ws.Range("K40").Value = 7000
$ws.Range("M40").Value = -6864
$ws.Range("H68").Value = 1112.7091
$ws.Range("I68").Value = 984.9808
$ws.Range("K68").Value = 984.9808
$ws.Range("M68").Value = -235.9808
$ws.Range("H71").Value = 1112.7091
$ws.Range("I71").Value = 984.9808
$ws.Range("K71").Value = 4924.904
$ws.Range("M71").Value = -1180.904
$ws.Range("H126").Value = 5642.143
$ws.Range("I126").Value = 2098.3333
$ws.Range("J126").Value = 8300
$ws.Range("K126").Value = 6294.999899999999
$ws.Range("L126").Value = 24900
$ws.Range("M126").Value = -3824.999899999999
$ws.Range("N126").Value = -29840
$ws.Range("H136").Value = 4663.625
$ws.Range("I136").Value = 1801.8
$ws.Range("J136").Value = 9433.333000000001
$ws.Range("K136").Value = 5405.4
$ws.Range("L136").Value = 28299.999
$ws.Range("M136").Value = -2855.4
$ws.Range("N136").Value = -33399.999

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12136.429
$ws.Range("I122").Value = 9500
$ws.Range("J122").Value = 13191
$ws.Range("K122").Value = 28500
$ws.Range("L122").Value = 39573
$ws.Range("M122").Value = -26050
$ws.Range("N122").Value = -44473
$ws.Range("H124").Value = 50429
$ws.Range("J124").Value = 50429
$ws.Range("L124").Value = 50429
$ws.Range("N124").Value = -60249
$ws.Range("H126").Value = 3134.8948
$ws.Range("I126").Value = 2111.25
$ws.Range("J126").Value = 3879.3635
$ws.Range("K126").Value = 6333.75
$ws.Range("L126").Value = 11638.0905
$ws.Range("M126").Value = -3863.75
$ws.Range("N126").Value = -16578.0905
$ws.Range("H132").Value = 12826435
$ws.Range("I132").Value = 11109.6
$ws.Range("J132").Value = 20836014
$ws.Range("K132").Value = 33328.8
$ws.Range("L132").Value = 62508042
$ws.Range("M132").Value = -30798.8
$ws.Range("H136").Value = 4283.4463
$ws.Range("I136").Value = 3627.3901
$ws.Range("J136").Value = 6076.6665
$ws.Range("K136").Value = 10882.1703
$ws.Range("L136").Value = 18229.9995
$ws.Range("M136").Value = -8332.1703
$ws.Range("N136").Value = -23329.9995
